$wb = $excel.ActiveWorkbook

# --- Fluxes sheet: remove the F09 flux (HerbZooplankton1 -> OmniZooplankton) ---
$fluxes = $wb.Worksheets.Item("Fluxes")
$fluxes.Rows.Item(6).Delete()

# Normalize the number format on the remaining boolean "Active" flux cells (D4:D5)
# so they share the same style as the rest of the sheet (drops the now-redundant
# duplicate style definition).
$fluxes.Range("D4").NumberFormat = "General"
$fluxes.Range("D5").NumberFormat = "General"

# --- Constraints sheet: add new biomass-range constraints C06-C09 ---
$constraints = $wb.Worksheets.Item("Constraints")

$constraints.Range("A7").Value = "C06"
$constraints.Range("B7").Value = "HerbZooplankton<=3*HerbZooplankton_Biomass"
$constraints.Range("C7").Value = "1988:1991"
$constraints.Range("D7").Value = 1

$constraints.Range("A8").Value = "C07"
$constraints.Range("B8").Value = "HerbZooplankton>=0.1*HerbZooplankton_Biomass"
$constraints.Range("C8").Value = "1988:1991"
$constraints.Range("D8").Value = 1

$constraints.Range("A9").Value = "C08"
$constraints.Range("B9").Value = "OmniZooplankton<=3*OmniZooplankton_Biomass"
$constraints.Range("C9").Value = "1988:1991"
$constraints.Range("D9").Value = 1

$constraints.Range("A10").Value = "C09"
$constraints.Range("B10").Value = "OmniZooplankton>=0.1*OmniZooplankton_Biomass"
$constraints.Range("C10").Value = "1988:1991"
$constraints.Range("D10").Value = 1
